$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-18 Wednesday" "2025-06-19 Thursday"
Replace-Text "95×34=3230" "17×56=952"
Replace-Text "41×45=1845" "98×86=8428"
Replace-Text "97×19=1843" "16×90=1440"
Replace-Text "63×92=5796" "65×17=1105"
Replace-Text "59×21=1239" "74×82=6068"
Replace-Text "66×28=1848" "69×89=6141"
Replace-Text "61×92=5612" "92×36=3312"
Replace-Text "50×92=4600" "72×71=5112"
Replace-Text "96×17=1632" "68×32=2176"
Replace-Text "83×65=5395" "45×82=3690"
Replace-Text "26×27=702" "22×24=528"
Replace-Text "72×79=5688" "81×84=6804"
Replace-Text "69×31=2139" "36×55=1980"
Replace-Text "39×67=2613" "48×40=1920"
Replace-Text "50×42=2100" "62×58=3596"
Replace-Text "45×77=3465" "47×35=1645"
Replace-Text "60×21=1260" "88×72=6336"
Replace-Text "23×18=414" "47×48=2256"
Replace-Text "31×39=1209" "49×86=4214"
Replace-Text "25×41=1025" "13×94=1222"
Replace-Text "71×59=4189" "62×60=3720"
Replace-Text "38×42=1596" "76×94=7144"
Replace-Text "94×53=4982" "86×99=8514"
Replace-Text "84×98=8232" "93×23=2139"
Replace-Text "95×23=2185" "70×95=6650"
